# Fehler in der Zeitberechnung wurden ausgebessert
# The "Planung" time value in B1 was wrong (60) and is corrected to 36.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B1").Value = 36
